$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix: the Springer API wasn't properly filling in the abstract in some
# situations. Re-fill the Abstract cell for row 4 with the full text, and
# re-write the Authors cells that had been affected by the same bug (their
# textual content is unchanged logically, but gets rewritten here).

$ws.Range("D4").Value = "`nMany European countries have responded to the COVID-19 pandemic by implementing nationwide protection measures and lockdowns1.`n However, the epidemic could rebound when such measures are relaxed, possibly leading to a requirement for a second or more, repeated lockdowns2.`n Here, we present results of a stochastic agent-based microsimulation model of the COVID-19 epidemic in France.`n We examined the potential impact of post-lockdown measures, including physical distancing, mask-wearing and shielding individuals who are the most vulnerable to severe COVID-19 infection, on cumulative disease incidence and mortality, and on intensive care unit (ICU)-bed occupancy.`n While lockdown is effective in containing the viral spread, once lifted, regardless of duration, it would be unlikely to prevent a rebound.`n Both physical distancing and mask-wearing, although effective in slowing the epidemic and in reducing mortality, would also be ineffective in ultimately preventing ICUs from becoming overwhelmed and a subsequent second lockdown.`n However, these measures coupled with the shielding of vulnerable people would be associated with better outcomes, including lower mortality and maintaining an adequate ICU capacity to prevent a second lockdown.`n Benefits would nonetheless be markedly reduced if most people do not adhere to these measures, or if they are not maintained for a sufficiently long period.`nA new study that models the COVID-19 epidemic in France suggests that although a second peak is likely unavoidable, maintaining social distancing and wearing masks when lockdown restrictions are lifted, as well as continuing to shelter vulnerable individuals, will reduce mortality and avoid overwhelming ICU facilities.`n"

$ws.Range("E2").Value = "[Jean%Dolbeault%xref no email%1,    Gabriel%Turinici%xref no email%1]"
$ws.Range("E4").Value = "[ Nicolas%Hoertel%null%0,  Martin%Blachier%null%1,  Carlos%Blanco%null%1,  Mark%Olfson%null%1,  Marc%Massetti%null%1,  Marina Sánchez%Rico%null%1,  Frédéric%Limosin%null%1,  Henri%Leleu%null%1]"
$ws.Range("E5").Value = "[Tuo%Ji%NULL%1,         Hai-Lian%Chen%NULL%1,         Jing%Xu%NULL%1,         Ling-Ning%Wu%NULL%1,         Jie-Jia%Li%NULL%1,         Kai%Chen%NULL%1,         Gang%Qin%tonygqin@ntu.edu.cn%1]"
$ws.Range("E6").Value = "[Can%Hou%xref no email%1,    Jiaxin%Chen%xref no email%1,    Yaqing%Zhou%xref no email%1,    Lei%Hua%xref no email%1,    Jinxia%Yuan%xref no email%1,    Shu%He%xref no email%1,    Yi%Guo%xref no email%1,    Sheng%Zhang%xref no email%2,    Qiaowei%Jia%xref no email%1,    Chenhui%Zhao%xref no email%1,    Jing%Zhang%xref no email%1,    Guangxu%Xu%xref no email%1,    Enzhi%Jia%xref no email%1]"
$ws.Range("E7").Value = "[Hien%Lau%hlau2@uci.edu%1,         Veria%Khosrawipour%veriakhosrawipour@yahoo.de%1,         Piotr%Kocbach%piotr.kocbach@uwm.edu.pl%1,         Agata%Mikolajczyk%agata.mikolajczyk@upwr.edu.pl%1,         Justyna%Schubert%justyna.schubert@upwr.edu.pl%1,         Jacek%Bania%jacek.bania@upwr.edu.pl%1,         Tanja%Khosrawipour%tkhosrawipour@gmail.com%1]"
$ws.Range("E8").Value = "[Carlo%Signorelli%NULL%1,         Thea%Scognamiglio%NULL%1,         Anna%Odone%NULL%1]"
$ws.Range("E9").Value = "[Jussi%Taipale%xref no email%1,    Paul%Romer%xref no email%1,    Sten%Linnarsson%xref no email%1]"
$ws.Range("E11").Value = "[Nadya%Johanna%NULL%1,         Henrico%Citrawijaya%NULL%1,         Grace%Wangge%NULL%1]"
$ws.Range("E12").Value = "[Gerard J.%Tellis%xref no email%1,    Ashish%Sood%xref no email%1,    Nitish%Sood%xref no email%1]"
$ws.Range("E13").Value = "[Biao%Tang%NULL%1,         Fan%Xia%NULL%1,         Sanyi%Tang%NULL%1,         Nicola Luigi%Bragazzi%NULL%1,         Qian%Li%NULL%1,         Xiaodan%Sun%NULL%1,         Juhua%Liang%NULL%1,         Yanni%Xiao%yxiao@mail.xjtu.edu.cn%1,         Jianhong%Wu%wujh@yorku.ca%1]"
